# Book Outline (Version 8) -- "Post Ch 8 1st draft" edit
#
# Applies the semantic/content changes described by the target diff:
#   1. "Examining Application and Services Logs" -> "Examining Event Logs"
#   2. "Managing PowerShell Script block logging" -> "Using PowerShell Script block logging"
#   3. Append a hyperlink (MSN article) right after "Managing Windows Defender"
#   4. Remove the stray empty paragraph just after "Installing WSL and WSL 2 (not needed)."
#   5. Remove the two stray empty ListParagraph paragraphs just after the
#      "New Material - 100%" that precedes the "9. Managing Storage" heading
#
# (Cosmetic/engine-generated artifacts from the original authoring session --
#  w:lastRenderedPageBreak relocation, w:proofErr spans, rsid/relationship-id
#  bookkeeping, run-splitting from literal keystrokes -- are not reproducible
#  through the Word object model and are intentionally left alone.)

$d = $word.ActiveDocument

function Get-ParaText($para) {
    return $para.Range.Text.TrimEnd([char]13)
}

function Find-ParagraphIndex($matchText, $occurrence) {
    $doc = $word.ActiveDocument
    $count = 0
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $t = Get-ParaText($doc.Paragraphs($i))
        if ($t -eq $matchText) {
            $count = $count + 1
            if ($count -eq $occurrence) {
                return $i
            }
        }
    }
    return -1
}

# ---------------------------------------------------------------------------
# 1. "Examining Application and Services Logs" -> "Examining Event Logs"
# ---------------------------------------------------------------------------
$idxExamining = Find-ParagraphIndex "Examining Application and Services Logs" 1
if ($idxExamining -gt 0) {
    $p = $d.Paragraphs($idxExamining)
    $rng = $d.Range($p.Range.Start, $p.Range.End - 1)
    $rng.Text = "Examining Event Logs"
} else {
    Write-Host "WARN: could not find 'Examining Application and Services Logs'"
}

# ---------------------------------------------------------------------------
# 2. "Managing PowerShell Script block logging" -> "Using PowerShell Script
#    block logging"
# ---------------------------------------------------------------------------
$idxPowerShell = Find-ParagraphIndex "Managing PowerShell Script block logging" 1
if ($idxPowerShell -gt 0) {
    $p = $d.Paragraphs($idxPowerShell)
    $rng = $d.Range($p.Range.Start, $p.Range.End - 1)
    $rng.Text = "Using PowerShell Script block logging"
} else {
    Write-Host "WARN: could not find 'Managing PowerShell Script block logging'"
}

# ---------------------------------------------------------------------------
# 3. Append hyperlink after "Managing Windows Defender"
# ---------------------------------------------------------------------------
$idxDefender = Find-ParagraphIndex "Managing Windows Defender" 1
if ($idxDefender -gt 0) {
    $url = "https://www.msn.com/en-gb/entertainment/news/lost-19th-century-tlingit-fort-discovered-in-alaska/ar-BB1d4Rnk"
    $p = $d.Paragraphs($idxDefender)
    $insertPos = $p.Range.End - 1
    $p.Range.InsertAfter($url)
    $afterEnd = $d.Paragraphs($idxDefender).Range.End - 1
    $hlRange = $d.Range($insertPos, $afterEnd)
    $d.Hyperlinks.Add($hlRange, $url)
} else {
    Write-Host "WARN: could not find 'Managing Windows Defender'"
}

# ---------------------------------------------------------------------------
# 4. Remove the stray empty paragraph after "Installing WSL and WSL 2 (not
#    needed)."
# ---------------------------------------------------------------------------
$idxWsl = Find-ParagraphIndex "Installing WSL and WSL 2 (not needed)." 1
if ($idxWsl -gt 0) {
    $stray = $d.Paragraphs($idxWsl + 1)
    if ((Get-ParaText $stray) -eq "") {
        $stray.Range.Delete()
    } else {
        Write-Host "WARN: paragraph after WSL line was not empty:" (Get-ParaText $stray)
    }
} else {
    Write-Host "WARN: could not find 'Installing WSL and WSL 2 (not needed).'"
}

# ---------------------------------------------------------------------------
# 5. Remove the two stray empty ListParagraph paragraphs that sit between
#    the second "New Material - 100%" and the "9. Managing Storage" heading.
# ---------------------------------------------------------------------------
$idxNewMaterial = Find-ParagraphIndex "New Material - 100%" 2
if ($idxNewMaterial -gt 0) {
    $n1 = $d.Paragraphs($idxNewMaterial + 1)
    $n2 = $d.Paragraphs($idxNewMaterial + 2)
    $n3 = $d.Paragraphs($idxNewMaterial + 3)
    if ((Get-ParaText $n1) -eq "" -and (Get-ParaText $n2) -eq "" -and (Get-ParaText $n3) -like "*Managing Storage*") {
        $d.Paragraphs($idxNewMaterial + 1).Range.Delete()
        $d.Paragraphs($idxNewMaterial + 1).Range.Delete()
    } else {
        Write-Host "WARN: unexpected context around 'New Material - 100%' (#2)"
    }
} else {
    Write-Host "WARN: could not find second 'New Material - 100%'"
}

Write-Host "Edits applied."
